$d = $word.ActiveDocument
$rng = $d.Content
$rng.Collapse(0)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p/><w:p/><w:p/><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="4944" w:type="pct"/><w:tblLayout w:type="fixed"/><w:tblLook w:val="04A0"/></w:tblPr><w:tblGrid><w:gridCol w:w="3708"/><w:gridCol w:w="5761"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/><w:shd w:val="clear" w:color="auto" w:fill="D9D9D9" w:themeFill="background1" w:themeFillShade="D9"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Property</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/><w:shd w:val="clear" w:color="auto" w:fill="D9D9D9" w:themeFill="background1" w:themeFillShade="D9"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Changes</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>proxy.siloDataAcessHost</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p/></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>proxy.siloDataAccessPort</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p/></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>mapserver.geonames.url</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t>http\://testteen.iwiglobal.com\:8081/geonames/servlet/iwiglobal?srv\=findNearbyAddress</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:r><w:t>wms.url</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t>https:://office.iwiglobal.com:8946/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cgi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-bin/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mapserv</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>wms.query</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t>?map=/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>var</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/www/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mapserver</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>romania</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/romania.map</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>wms.layers</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>romania</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Rom </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rompoly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> POI</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>wms.layer.queryparam</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>rompoi</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>google.map.url</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p/></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>google.map.geo.url</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p/></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>mailserver</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t>testteen.iwiglobal.com</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>mail.debug</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t>false</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>teenproxy.siloDataAccessHost</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t>10.0.35.20</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>teenproxy.silodataaccessport</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t>8099</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>teenproxy.centraldataaccesshost</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t>10.0.35.20</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>teenproxy.centraldataaccessport</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t>8095</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>teenproxy.sbsUserId</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t>527</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>teenproxy.sbsUserName</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>tiwipro</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>speedbystreet.iwigis</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t>192.168.11.80</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>speedbystreet.username</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t>development</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>speedbystreet.password</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t>w1tn3ss</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>cron.emailreport</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t xml:space="preserve">0 1/15 * * </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>* ?</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>cron.emailalert</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t xml:space="preserve">0 1/2 * * </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>* ?</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>cronsmsalert</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t xml:space="preserve">0 0/2 * * </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>* ?</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>cron.phonealert</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t xml:space="preserve">30 1/2 * * </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>* ?</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>accountOptions.phoneAlerts</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t>MyAccountBean.java, PersonBean.java</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>accountOptions.noReplyEmailAddress</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t>applicationContext-reportx.xml</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:r><w:t>externalConfig.phone01</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t xml:space="preserve">Used for error pages (403.xhtml, 404.xhtml,500.xhtml) and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>inthincCustomerSupport.xhtml</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:r><w:t>externalConfig.phone02</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:r><w:t xml:space="preserve">Used for error pages (403.xhtml, 404.xhtml,500.xhtml) and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>inthincCustomerSupport.xhtml</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1958" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>addressLookup</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3042" w:type="pct"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>googleAjaxAddressLookupBean</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc></w:tr></w:tbl></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)
